$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "Total nominations received this Session" row (row 38).
# Everything below shifts up by one; the row 36/37/38-41 relabeling below
# then lines the text up with the new totals layout.
$ws.Rows("38").Delete()

# --- Section label prefixing (A7:A35) ---
# Civilian section
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

# Other Civilian section
$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Confirmed "
$ws.Range("A15").Value = "     Other Civilian, Unconfirmed "
$ws.Range("A16").Value = "     Other Civilian, Returned to White House "

# Air Force section
$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("A19").Value = "     Air Force, Confirmed "
$ws.Range("A20").Value = "     Air Force, Unconfirmed "
$ws.Range("A21").Value = "     Air Force, Withdrawn "
$ws.Range("A22").Value = "     Air Force, Returned to White House "

# Army section
$ws.Range("A24").Value = "     Army, New nominations"
$ws.Range("A25").Value = "     Army, Confirmed "
$ws.Range("A26").Value = "     Army, Unconfirmed "
$ws.Range("A27").Value = "     Army, Withdrawn "

# Navy section
$ws.Range("A29").Value = "     Navy, New nominations"
$ws.Range("A30").Value = "     Navy, Confirmed "
$ws.Range("A31").Value = "     Navy, Unconfirmed "
$ws.Range("A32").Value = "     Navy, Returned to White House "

# Marine Corps section
$ws.Range("A34").Value = "     Marine Corps, New nominations"
$ws.Range("A35").Value = "     Marine Corps, Confirmed "

# --- Summary section (rows 36-41 after the row-38 delete above) ---
$ws.Range("A36").Value = "Total new nominations"
$ws.Range("B36").Value = 19999
$ws.Range("B36").NumberFormat = "#,##0"

$ws.Range("A37").Value = "Total carryover nominations"

$ws.Range("A38").Value = "Total confirmed "
$ws.Range("A39").Value = "Total unconfirmed "
$ws.Range("A40").Value = "Total withdrawn "
$ws.Range("A41").Value = "Total returned to the White House "
